$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 716; this shifts existing rows 716-748 down to 717-749
$ws.Rows(716).Insert()

# Populate the newly inserted row 716 with the new weekly price record
$ws.Cells.Item(716, 1).Value = 4
$ws.Cells.Item(716, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(716, 3).Value = "Los Lagos"
$ws.Cells.Item(716, 4).Value = 45041
$ws.Cells.Item(716, 5).Value = 10
$ws.Cells.Item(716, 6).Value = "Fruta"
$ws.Cells.Item(716, 7).Value = 100102
$ws.Cells.Item(716, 8).Value = "Cítricos"
$ws.Cells.Item(716, 9).Value = 100102005
$ws.Cells.Item(716, 10).Value = "Naranja"
$ws.Cells.Item(716, 11).Value = "Valencia"
$ws.Cells.Item(716, 12).Value = "Primera"
$ws.Cells.Item(716, 13).Value = 600
$ws.Cells.Item(716, 14).Value = 19000
$ws.Cells.Item(716, 15).Value = 20000
$ws.Cells.Item(716, 16).Value = 19500
$ws.Cells.Item(716, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(716, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(716, 19).Value = 1300
$ws.Cells.Item(716, 20).Value = 15
